$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.363.98'
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").Value = '1.933.17'
$ws.Range("E3").Value = '  +0.72%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.54%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.95'
$ws.Range("E5").Value = '  +2.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7128'
$ws.Range("E6").Value = '  -1.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.003'
$ws.Range("E7").Value = '  +0.60%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3259'
$ws.Range("E8").Value = '  +0.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '27.52'
$ws.Range("E9").Value = '  +4.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07159'
$ws.Range("E10").Value = '  +4.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7999'
$ws.Range("E11").Value = '  +1.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08082'
$ws.Range("E12").Value = '  +2.10%  '
$ws.Range("D13").Value = '1.932.28'
$ws.Range("E13").Value = '  +0.76%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.416'
$ws.Range("E14").Value = '  +0.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.44'
$ws.Range("E15").Value = '  +0.42%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.91'
$ws.Range("E16").Value = '  +3.42%  '
$ws.Range("D17").Value = '30.340.25'
$ws.Range("E17").Value = '  +0.44%  '
$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '252.20'
$ws.Range("E18").Value = '  -2.64%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008090'
$ws.Range("E19").Value = '  +2.97%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.800'
$ws.Range("E20").Value = '  -0.21%  '
$ws.Range("D21").Value = '2.185.50'
$ws.Range("E21").Value = '  +0.97%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.004'
$ws.Range("E23").Value = '  +0.66%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.925'
$ws.Range("E24").Value = '  +1.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.683'
$ws.Range("E25").Value = '  +0.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.50'
$ws.Range("E26").Value = '  +3.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.331'
$ws.Range("E27").Value = '  +5.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.19'
$ws.Range("E28").Value = '  +2.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1284'
$ws.Range("E29").Value = '  -3.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.367'
$ws.Range("E30").Value = '  +0.87%  '
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.422'
$ws.Range("E32").Value = '  +0.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.194'
$ws.Range("E33").Value = '  +0.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05204'
$ws.Range("E34").Value = '  +3.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.272'
$ws.Range("E35").Value = '  +7.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7468'
$ws.Range("E36").Value = '  +1.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.765'
$ws.Range("E37").Value = '  +1.85%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01964'
$ws.Range("E38").Value = '  +1.57%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.804'
$ws.Range("E39").Value = '  +0.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '78.96'
$ws.Range("E40").Value = '  -0.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.455'
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4522'
$ws.Range("E42").Value = '  +2.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.022'
$ws.Range("E43").Value = '  +0.58%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.002'
$ws.Range("E44").Value = '  +0.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8407'
$ws.Range("E45").Value = '  +1.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.72'
$ws.Range("E46").Value = '  -0.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.768'
$ws.Range("E47").Value = '  +1.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.406'
$ws.Range("E48").Value = '  +2.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.64'
$ws.Range("E49").Value = '  +1.82%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4174'
$ws.Range("E50").Value = '  +2.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06060'
$ws.Range("E51").Value = '  +2.60%  '
